# Update "want to go" (想去人数) counts in column F across all sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 322
$ws.Range("F6").Value = 400
$ws.Range("F7").Value = 872
$ws.Range("F8").Value = 55
$ws.Range("F9").Value = 520
$ws.Range("F10").Value = 68
$ws.Range("F12").Value = 1147
$ws.Range("F14").Value = 244
$ws.Range("F15").Value = 36
$ws.Range("F16").Value = 417
$ws.Range("F17").Value = 6658
$ws.Range("F18").Value = 65
$ws.Range("F21").Value = 7589
$ws.Range("F24").Value = 3403
$ws.Range("F25").Value = 29
$ws.Range("F26").Value = 2110
$ws.Range("F27").Value = 896
$ws.Range("F29").Value = 145
$ws.Range("F32").Value = 225
$ws.Range("F33").Value = 199
$ws.Range("F34").Value = 1705
$ws.Range("F36").Value = 175
$ws.Range("F39").Value = 1213
$ws.Range("F40").Value = 1793
$ws.Range("F41").Value = 2140

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 20
$ws.Range("F7").Value = 82

# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 1231
$ws.Range("F4").Value = 77

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1231
$ws.Range("F5").Value = 77
$ws.Range("F7").Value = 322
$ws.Range("F8").Value = 400
$ws.Range("F9").Value = 872
$ws.Range("F10").Value = 55
$ws.Range("F11").Value = 520
$ws.Range("F14").Value = 1147
$ws.Range("F17").Value = 244
$ws.Range("F18").Value = 36
$ws.Range("F19").Value = 417
$ws.Range("F20").Value = 6658
$ws.Range("F21").Value = 65
$ws.Range("F24").Value = 7589
$ws.Range("F27").Value = 3403
$ws.Range("F28").Value = 29
$ws.Range("F29").Value = 2110
$ws.Range("F30").Value = 896
$ws.Range("F32").Value = 145
$ws.Range("F36").Value = 225
$ws.Range("F37").Value = 199
$ws.Range("F38").Value = 1705
$ws.Range("F40").Value = 175
$ws.Range("F44").Value = 1213
$ws.Range("F45").Value = 1793
$ws.Range("F46").Value = 20
$ws.Range("F47").Value = 2140
$ws.Range("F49").Value = 82
